$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 34; $row++) {
    $cell = $ws.Range("M$row")
    if ($cell.Value2 -eq "Catolicismo") {
        $cell.Value = "Catolica"
    }
}
